$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-4 with the new data
$ws.Range("A2").Value = "#borabill"
$ws.Range("B2").Value = "Bora Biil"
$ws.Range("C2").Value = 4000
$ws.Range("D2").Value = "Dança da música Bora Bill"

$ws.Range("A3").Value = "#elizabeth"
$ws.Range("B3").Value = "Pensava que era imortal"
$ws.Range("C3").Value = 7000
$ws.Range("D3").Value = "Vídeo com narração chorando por que pensava que a Rainha Elizabeth era imortal"

$ws.Range("A4").Value = "#sonho"
$ws.Range("B4").Value = "Pensava que era realidad3"
$ws.Range("C4").Value = 9000
$ws.Range("D4").Value = "Finge que ta sonhando que ganhou milhões de reais e acorda com gritos da mãe"

# Remove the old rows 5 and 6 entirely (shrinks the used range back to A1:D4)
$ws.Range("A5:D6").Delete()

# Update the selected cell to match the new last row (D4)
$ws.Range("D4").Select() | Out-Null
